# Regenerate orders with updated distance/size codes.
#
# The authored change renames:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31
# wherever these substrings occur inside text values on the sheet
# (Condition / Filename_Left / Filename_Right / Distance / Size columns,
# including composite strings like "Face15_D64_S25" and
# "Face15_D64_S25_l.png"). Numeric / boolean cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val
            $newVal = $newVal -replace 'D64', 'D69'
            $newVal = $newVal -replace 'D51', 'D55'
            $newVal = $newVal -replace 'D80', 'D86'
            $newVal = $newVal -replace 'S30', 'S31'

            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
